# Apply the "Added and adapted Yarpiz PSO implementation" edit.
#
# Summary of the required change:
#  - A1 header changes from "Gen" to "MaxFES"
#  - A2:A14 values change from generation counts to the fraction-of-budget
#    values 0, 0.001, 0.01, 0.1, 0.2, ..., 1
#  - The "Run 50" column (AZ) is removed; the old "Mean" column (BA) is
#    recomputed (mean over only the remaining 50 runs, columns B:AY) and
#    becomes the new last column (AZ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header of column A.
$ws.Range("A1").Value = "MaxFES"

# 2. Update column A's data values (rows 2-14).
$maxfesValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $maxfesValues.Length; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $maxfesValues[$i]
}

# 3. Recompute the "Mean" column (currently BA) so that it only averages
#    the 50 remaining runs (columns B:AY), excluding the about-to-be-removed
#    "Run 50" column (AZ).
for ($row = 2; $row -le 14; $row++) {
    $ws.Range("BA$row").Formula = "=AVERAGE(B$row" + ":AY$row)"
}

# Convert the freshly computed formulas to plain static values, matching
# the rest of the sheet (which stores literal numbers, not formulas).
for ($row = 2; $row -le 14; $row++) {
    $computed = $ws.Range("BA$row").Value2
    $ws.Range("BA$row").Value = $computed
}

# 4. Remove the "Run 50" column (AZ). This shifts the recomputed "Mean"
#    column (currently BA) left into AZ, becoming the new final column.
$ws.Columns("AZ").Delete()
